$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common values shared by the new rows (39-41 share most fields; 42-43 share most fields)
$companyName = "Western Interior Designers & Marine Contractors"
$date = "16-01-2026"
$corpId = 286962
$companyAcct = 34413429360
$txnType = "NEFT"
$companyIfsc = "SBIN0003229"
$companyPan = "AAAFW8862C"
$companyGstin = "32AAAFW8862C1Z9"
$status = "pending"
$fromMail = "hrm@westernidc.com"
$estStatus = "ESTIMATION NOT MATCHED"

# Row 39
$r = 39
$ws.Cells.Item($r, 1).Value = "WGG 02"
$ws.Cells.Item($r, 2).Value = $companyName
$ws.Cells.Item($r, 3).Value = $date
$ws.Cells.Item($r, 4).Value = $corpId
$ws.Cells.Item($r, 5).Value = $companyName
$ws.Cells.Item($r, 6).Value = $companyAcct
$ws.Cells.Item($r, 7).Value = $txnType
$ws.Cells.Item($r, 8).Value = $companyIfsc
$ws.Cells.Item($r, 9).Value = $companyPan
$ws.Cells.Item($r, 10).Value = $companyGstin
$ws.Cells.Item($r, 12).Value = "6186c0de-5d0d-4300-8f62-ddaf3e547736"
$ws.Cells.Item($r, 21).Value = $status
$ws.Cells.Item($r, 22).Value = 9900
$ws.Cells.Item($r, 24).Value = "DUMKA ROOM RENT RPA_UNIQUE_ID : e8e519e2-63ea-40a3-9e23-82d179093abb"
$ws.Cells.Item($r, 25).Value = "dumka"
$ws.Cells.Item($r, 26).Value = 0
$ws.Cells.Item($r, 27).Value = $fromMail
$ws.Cells.Item($r, 28).Value = $estStatus
$ws.Cells.Item($r, 29).Value = 0
$ws.Cells.Item($r, 30).Value = 0
$ws.Cells.Item($r, 31).Value = 0

# Row 40
$r = 40
$ws.Cells.Item($r, 1).Value = "WGG 02"
$ws.Cells.Item($r, 2).Value = $companyName
$ws.Cells.Item($r, 3).Value = $date
$ws.Cells.Item($r, 4).Value = $corpId
$ws.Cells.Item($r, 5).Value = $companyName
$ws.Cells.Item($r, 6).Value = $companyAcct
$ws.Cells.Item($r, 7).Value = $txnType
$ws.Cells.Item($r, 8).Value = $companyIfsc
$ws.Cells.Item($r, 9).Value = $companyPan
$ws.Cells.Item($r, 10).Value = $companyGstin
$ws.Cells.Item($r, 12).Value = "d99570da-b209-4da6-ade8-4e678cbd864c"
$ws.Cells.Item($r, 21).Value = $status
$ws.Cells.Item($r, 22).Value = 4000
$ws.Cells.Item($r, 24).Value = "COOK SALARY DUMKA RPA_UNIQUE_ID : 164d8204-fc59-4968-b6ad-6adc54bd815a"
$ws.Cells.Item($r, 25).Value = "dumka"
$ws.Cells.Item($r, 26).Value = 0
$ws.Cells.Item($r, 27).Value = $fromMail
$ws.Cells.Item($r, 28).Value = $estStatus
$ws.Cells.Item($r, 29).Value = 0
$ws.Cells.Item($r, 30).Value = 0
$ws.Cells.Item($r, 31).Value = 0

# Row 41
$r = 41
$ws.Cells.Item($r, 1).Value = "WGG 02"
$ws.Cells.Item($r, 2).Value = $companyName
$ws.Cells.Item($r, 3).Value = $date
$ws.Cells.Item($r, 4).Value = $corpId
$ws.Cells.Item($r, 5).Value = $companyName
$ws.Cells.Item($r, 6).Value = $companyAcct
$ws.Cells.Item($r, 7).Value = $txnType
$ws.Cells.Item($r, 8).Value = $companyIfsc
$ws.Cells.Item($r, 9).Value = $companyPan
$ws.Cells.Item($r, 10).Value = $companyGstin
$ws.Cells.Item($r, 12).Value = "c3226756-1dbc-46d3-b099-af90812ae3b3"
$ws.Cells.Item($r, 21).Value = $status
$ws.Cells.Item($r, 22).Value = 3420
$ws.Cells.Item($r, 24).Value = "GROCERY PURCHASE EXPENSES RPA_UNIQUE_ID : e01c71fb-5a74-44ec-bac6-db672e455f97"
$ws.Cells.Item($r, 25).Value = "dumka"
$ws.Cells.Item($r, 26).Value = 0
$ws.Cells.Item($r, 27).Value = $fromMail
$ws.Cells.Item($r, 28).Value = $estStatus
$ws.Cells.Item($r, 29).Value = 0
$ws.Cells.Item($r, 30).Value = 0
$ws.Cells.Item($r, 31).Value = 0

# Row 42
$r = 42
$ws.Cells.Item($r, 1).Value = "WGE 77"
$ws.Cells.Item($r, 2).Value = $companyName
$ws.Cells.Item($r, 3).Value = $date
$ws.Cells.Item($r, 4).Value = $corpId
$ws.Cells.Item($r, 5).Value = $companyName
$ws.Cells.Item($r, 6).Value = $companyAcct
$ws.Cells.Item($r, 7).Value = $txnType
$ws.Cells.Item($r, 8).Value = $companyIfsc
$ws.Cells.Item($r, 9).Value = $companyPan
$ws.Cells.Item($r, 10).Value = $companyGstin
$ws.Cells.Item($r, 11).Value = "Anju M S"
$ws.Cells.Item($r, 12).Value = "566dfec3-c7ce-42d3-b94f-a1528fb41bcb"
$ws.Cells.Item($r, 13).Value = 345002010013320
$ws.Cells.Item($r, 14).Value = "UBIN0534501"
$ws.Cells.Item($r, 21).Value = $status
$ws.Cells.Item($r, 22).Value = 4300
$ws.Cells.Item($r, 24).Value = "material shifting payment to mdl to ipshem RPA_UNIQUE_ID : 1d2a7ae7-44da-4597-a109-c9d8089a7576"
$ws.Cells.Item($r, 25).Value = "mdl mumbai"
$ws.Cells.Item($r, 26).Value = 0
$ws.Cells.Item($r, 27).Value = $fromMail
$ws.Cells.Item($r, 28).Value = $estStatus
$ws.Cells.Item($r, 29).Value = 0
$ws.Cells.Item($r, 30).Value = 0
$ws.Cells.Item($r, 31).Value = 0

# Row 43
$r = 43
$ws.Cells.Item($r, 1).Value = "WGE 77"
$ws.Cells.Item($r, 2).Value = $companyName
$ws.Cells.Item($r, 3).Value = $date
$ws.Cells.Item($r, 4).Value = $corpId
$ws.Cells.Item($r, 5).Value = $companyName
$ws.Cells.Item($r, 6).Value = $companyAcct
$ws.Cells.Item($r, 7).Value = $txnType
$ws.Cells.Item($r, 8).Value = $companyIfsc
$ws.Cells.Item($r, 9).Value = $companyPan
$ws.Cells.Item($r, 10).Value = $companyGstin
$ws.Cells.Item($r, 11).Value = "Anju M S"
$ws.Cells.Item($r, 12).Value = "2ecd8d46-f937-4c21-a7e8-c8cbc46c196e"
$ws.Cells.Item($r, 13).Value = 345002010013320
$ws.Cells.Item($r, 14).Value = "UBIN0534501"
$ws.Cells.Item($r, 21).Value = $status
$ws.Cells.Item($r, 22).Value = 6600
$ws.Cells.Item($r, 24).Value = "DUMKA ROOM RENT ( 02/01/2026 to 07/01/2026) This amount has already paid by anju/ so the  amount credited to anjus account RPA_UNIQUE_ID : 40f014e1-9fb2-4300-9f38-6428269a7fe4"
$ws.Cells.Item($r, 25).Value = "dumka"
$ws.Cells.Item($r, 26).Value = 0
$ws.Cells.Item($r, 27).Value = $fromMail
$ws.Cells.Item($r, 28).Value = $estStatus
$ws.Cells.Item($r, 29).Value = 0
$ws.Cells.Item($r, 30).Value = 0
$ws.Cells.Item($r, 31).Value = 0
